$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new client entry "CSC" in column A, row 3 (Region column)
$ws.Range("A3").Value = "CSC"

# Update the selected cell/range as recorded in the saved view
$ws.Range("D10").Select()
